$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.057.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.63"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6319"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07540"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.21%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.13"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07706"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.99"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.000"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6701"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.19"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009589"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.080"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.071.38"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.58"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "226.58"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.142"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.95%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1429"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.517"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.96"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.150"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.066"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05480"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.202"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.859"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7448"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.140"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.655"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.245.96"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.753"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.601"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9034"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.39"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.978.35"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.04"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000122"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.79%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.016"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.658"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.89%  "
